# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 85, shifting the existing
# historical rows (old 85-160) down to 86-161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 85 — this pushes every
# row from 85 downward (old row85 -> new row86, ... old row160 -> new row161)
# and also shifts the sheet's dimension/used range accordingly.
$ws.Rows("85").Insert()

# Populate the freshly inserted row 85 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T keep the same values used throughout
# this data block; D,M,N,O,P,R,S carry the new record's data.
$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value = 45280
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100101
$ws.Range("H85").Value = "Berries"
$ws.Range("I85").Value = 100101001
$ws.Range("J85").Value = "Arándano (blue)"
$ws.Range("K85").Value = "Sin especificar"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 180
$ws.Range("N85").Value = 4000
$ws.Range("O85").Value = 4500
$ws.Range("P85").Value = 4278
$ws.Range("Q85").Value = "$/bandeja 2 kilos"
$ws.Range("R85").Value = "Región de Ñuble"
$ws.Range("S85").Value = 2139
$ws.Range("T85").Value = 2
